$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: set C13 and D13 to 5
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 5

# Row 17: set C17, D17, E17 to 5
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 5

# Row 19: set D19 and E19 to 5 (C19 already 5)
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 5

# Update the active selection to F19
$ws.Range("F19").Select()
